$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: was a blank separator row, now becomes the first bolt line ---
$ws.Range("B31").Value = "M2-12mm Bolt"
$ws.Range("C31").Value = 1
$ws.Range("E31").Value = "any store"

# --- Rows 32-41: existing bolt/nut rows, values updated in place ---
$ws.Range("B32").Value = "M3-12mm Bolt"
$ws.Range("C32").Value = 44

$ws.Range("B33").Value = "M3-5mm Bolt"
$ws.Range("C33").Value = 4

$ws.Range("B34").Value = "M3-25mm Bolt"
$ws.Range("C34").Value = 2

$ws.Range("B35").Value = "M3-20mm Bolt"
$ws.Range("C35").Value = 6

$ws.Range("B36").Value = "M3-40mm Bolt"
$ws.Range("C36").Value = 1

$ws.Range("B37").Value = "M3.5-10mm Bolt"
$ws.Range("C37").Value = 6

$ws.Range("B38").Value = "M3.5-35mm Bolt"
$ws.Range("C38").Value = 8

$ws.Range("B39").Value = "M5-20mm Bolt"
$ws.Range("C39").Value = 2

$ws.Range("B40").Value = "M2 nut"
$ws.Range("C40").Value = 1

$ws.Range("B41").Value = "M3 nut"
$ws.Range("C41").Value = 16

# --- Row 42: used to only carry "any store" in column E; fill B/C ---
$ws.Range("B42").Value = "M3.5 nut"
$ws.Range("C42").Value = 8

# --- Insert 3 fresh rows before the "Prints" header (old row 43) ---
$ws.Rows("43:45").Insert()
$ws.Rows("43:45").RowHeight = 30

$ws.Range("B43").Value = "M5 nut"
$ws.Range("C43").Value = 2
$ws.Range("E43").Value = "any store"

$ws.Range("B44").Value = "M3-10mm nylon Screw"
$ws.Range("C44").Value = 4
$ws.Range("E44").Value = "any store"

$ws.Range("B45").Value = "M3 nylon nut"
$ws.Range("C45").Value = 4
$ws.Range("E45").Value = "any store"

# --- Grow the CustomerList table to cover the three new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B4:F62"))

# --- Restore the selection to match the saved view ---
$ws.Activate()
$ws.Range("E35").Select()
